$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 10.73888328386894
    "C2" = 8.185796483737697
    "E2" = 12.07696957496907
    "F2" = 16.86991607391245
    "G2" = 3.595917274550334
    "I2" = 16.94453677806763
    "M2" = 13.93742979277052
    "N2" = 16.24575604260859
    "O2" = 17.65193695639881
    "B3" = 10.18977078708362
    "C3" = 7.795986312625922
    "E3" = 11.95620324643995
    "F3" = 15.89584955866815
    "G3" = 3.597831346401437
    "I3" = 17.04296267437787
    "M3" = 13.66519876601906
    "N3" = 16.29171703176359
    "O3" = 17.69623411241293
    "B4" = 9.837655821444349
    "C4" = 7.544848478787079
    "E4" = 11.88610524284747
    "F4" = 15.26997757108489
    "G4" = 3.599068650462961
    "I4" = 17.10775544134084
    "M4" = 13.49842400350246
    "N4" = 16.3217386887834
    "O4" = 17.72895354930095
    "B5" = 9.690562912939214
    "C5" = 7.439615747037998
    "E5" = 11.85859071970586
    "F5" = 15.008197319934
    "G5" = 3.599588516392122
    "I5" = 17.13525250457576
    "M5" = 13.43065306264174
    "N5" = 16.33442654168298
    "O5" = 17.74366893412532
    "B6" = 9.665925792904549
    "C6" = 7.421969661748911
    "E6" = 11.854086308848
    "F6" = 14.96433081551589
    "G6" = 3.599675786694422
    "I6" = 17.13988434040085
    "M6" = 13.41941414192501
    "N6" = 16.33656078128115
    "O6" = 17.74619568484928
    "B7" = 9.835686435006572
    "C7" = 7.543440876606029
    "E7" = 11.88572987690027
    "F7" = 15.26647399323133
    "G7" = 3.599075598105989
    "I7" = 17.1081218522671
    "M7" = 13.49750911843052
    "N7" = 16.32190796318586
    "O7" = 17.7291464193987
    "B8" = 10.55274025339237
    "C8" = 8.053877671184217
    "E8" = 12.0345080255022
    "F8" = 16.5399640634477
    "G8" = 3.596564396875527
    "I8" = 16.97756777318229
    "M8" = 13.84354386863648
    "N8" = 16.26122998501796
    "O8" = 17.6660611332669
    "B9" = 11.83454855048406
    "C9" = 8.958655810235994
    "E9" = 12.35689947845454
    "F9" = 19.00274580682531
    "G9" = 3.592130045971353
    "I9" = 16.75625477312876
    "M9" = 14.52095751544048
    "N9" = 16.15649845815053
    "O9" = 17.5864153633648
    "B10" = 12.6944540515292
    "C10" = 9.561975194110911
    "E10" = 12.61025272208073
    "F10" = 20.67494806633232
    "G10" = 3.589167714531277
    "I10" = 16.61497576221922
    "M10" = 15.01258886047631
    "N10" = 16.08819545042464
    "O10" = 17.55508224728538
    "B11" = 13.06696005830244
    "C11" = 9.822685844207601
    "E11" = 12.72859283938438
    "F11" = 21.3917225636224
    "G11" = 3.587883576663899
    "I11" = 16.55537042917009
    "M11" = 15.23386393827855
    "N11" = 16.0589894622626
    "O11" = 17.54678653307923
    "B12" = 13.20527584633053
    "C12" = 9.919407761494504
    "E12" = 12.7738061135652
    "F12" = 21.65686569030329
    "G12" = 3.587406378324558
    "I12" = 16.53347328202302
    "M12" = 15.31723125611906
    "N12" = 16.04819743931175
    "O12" = 17.54450543360879
    "B13" = 13.1756099388432
    "C13" = 9.898666451966472
    "E13" = 12.76405151700249
    "F13" = 21.60004134736742
    "G13" = 3.587508748560986
    "I13" = 16.53815918650622
    "M13" = 15.2992967885568
    "N13" = 16.05050979986624
    "O13" = 17.54495839723851
    "B14" = 13.07839467812786
    "C14" = 9.830683505875783
    "E14" = 12.73230478921712
    "F14" = 21.4136618050453
    "G14" = 3.587844135618944
    "I14" = 16.55355540719393
    "M14" = 15.24073159393054
    "N14" = 16.05809623522019
    "O14" = 17.54658160579817
    "B15" = 13.01848858652268
    "C15" = 9.788780349340804
    "E15" = 12.71290984574249
    "F15" = 21.29868154950795
    "G15" = 3.588050750583745
    "I15" = 16.56307393089424
    "M15" = 15.20480099738649
    "N15" = 16.06277798483173
    "O15" = 17.54768799961257
    "B16" = 12.66972887343195
    "C16" = 9.54465813491357
    "E16" = 12.60257734426995
    "F16" = 20.62722412089977
    "G16" = 3.589252908625676
    "I16" = 16.61896524677548
    "M16" = 14.99807303595482
    "N16" = 16.09014161594659
    "O16" = 17.55574460040702
    "B17" = 12.45094823049283
    "C17" = 9.391355968976152
    "E17" = 12.53565207606377
    "F17" = 20.20408069597325
    "G17" = 3.590006609699905
    "I17" = 16.65444960524914
    "M17" = 14.87058382596899
    "N17" = 16.10740566220278
    "O17" = 17.56221577728838
    "B18" = 12.32335687562668
    "C18" = 9.301888364059019
    "E18" = 12.49745179532529
    "F18" = 19.95656407809801
    "G18" = 3.590446092771318
    "I18" = 16.67529797506299
    "M18" = 14.79703701026978
    "N18" = 16.11751109843586
    "O18" = 17.56649843964932
    "B19" = 12.27985730568351
    "C19" = 9.271375153941126
    "E19" = 12.48456954312343
    "F19" = 19.87204792380568
    "G19" = 3.59059592166172
    "I19" = 16.68243211657206
    "M19" = 14.7721004453784
    "N19" = 16.1209628056473
    "O19" = 17.56804464130468
    "B20" = 12.47441985771054
    "C20" = 9.40780920614557
    "E20" = 12.54274633771339
    "F20" = 20.24955283636154
    "G20" = 3.589925758965716
    "I20" = 16.65062680498149
    "M20" = 14.88417849553163
    "N20" = 16.10554970400806
    "O20" = 17.56146886008924
    "B21" = 13.10702404920988
    "C21" = 9.850706311657962
    "E21" = 12.74161905547721
    "F21" = 21.46857628470577
    "G21" = 3.587745378333074
    "I21" = 16.54901484051337
    "M21" = 15.25794578381964
    "N21" = 16.05586065588498
    "O21" = 17.5460814565438
    "B22" = 13.50445262612414
    "C22" = 10.12847800425811
    "E22" = 12.87390740428401
    "F22" = 22.22866616901552
    "G22" = 3.586373258984517
    "I22" = 16.48653641245682
    "M22" = 15.49971301747989
    "N22" = 16.02494580645992
    "O22" = 17.54104059582058
    "B23" = 13.293819626311
    "C23" = 9.981303056013685
    "E23" = 12.80310547748401
    "F23" = 21.82633154458858
    "G23" = 3.587100760997775
    "I23" = 16.51952146070727
    "M23" = 15.37093355415427
    "N23" = 16.04130310952661
    "O23" = 17.54327106023528
    "B24" = 12.46381397266785
    "C24" = 9.400374844892776
    "E24" = 12.5395381591626
    "F24" = 20.22900810905287
    "G24" = 3.589962292358873
    "I24" = 16.65235369791096
    "M24" = 14.87803312256307
    "N24" = 16.10638822162515
    "O24" = 17.56180479008318
    "B25" = 11.50179437025569
    "C25" = 8.724509085361385
    "E25" = 12.26662275854332
    "F25" = 18.34778573295695
    "G25" = 3.593277517731137
    "I25" = 16.81239400211089
    "M25" = 14.33838944137122
    "N25" = 16.18330983130021
    "O25" = 17.60320891570439
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
